$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Row 10 (Objetivos:) - replace B/C with new objectives text
$ws.Range("B10").Value = "Propiciar ao aluno conhecimentos básicos da Química Inorgânica envolvida em processos biológicos."
$ws.Range("C10").Value = "Propiciar ao aluno conhecimentos básicos da Química Inorgânica envolvida em processos biológicos."

# 2) Insert new row at 13 for "2143261 - Andre Luis Ferraz" (Docentes responsaveis row)
$ws.Rows(13).Insert()
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "2143261 - André Luis Ferraz"
$ws.Range("C13").Value = "2143261 - André Luis Ferraz"
$ws.Range("B14:C14").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)

# 3) Row 14 (Programa resumido:) - replace B/C with new summary text
$ws.Range("B14").Value = "Estrutura molecular e ligação química; Orbitais moleculares e as moléculas de O2 e N2; Ácidos, bases e a correlação com os ligantes dos metais em solução; Complexos metálicos - teoria do campo cristalino; Sistemas biológicos de transporte; Transporte de O2 e transferência de elétrons em sistemas biológicos; Processos catalíticos - ácido/base e oxido-redução em metaloproteínas."
$ws.Range("C14").Value = "Estrutura molecular e ligação química; Orbitais moleculares e as moléculas de O2 e N2; Ácidos, bases e a correlação com os ligantes dos metais em solução; Complexos metálicos - teoria do campo cristalino; Sistemas biológicos de transporte; Transporte de O2 e transferência de elétrons em sistemas biológicos; Processos catalíticos - ácido/base e oxido-redução em metaloproteínas."

# 4) Row 16 (Programa:) - replace B/C with new detailed program text
$ws.Range("B16").Value = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."
$ws.Range("C16").Value = "1. Estrutura molecular e ligação química: Teoria de ligação de valência, estrutura de compostos com C, N, O; Relação entre estrutura e propriedades fisico-químicas2. Orbitais moleculares e as moléculas de O2 e N2: Limitações da teoria de ligação de valência, reatividade diferenciada de O2 e N2, relevância do O2 em sistemas biológicos, espécies reativas de oxigênio3. Ácidos, bases e a correlação com os ligantes dos metais em solução: Afinidade das bases por metais de transição, equilíbrio químico em sistemas biológicos4. Complexos metálicos - teoria do campo cristalino: Teoria do campo cristalino e os compostos octaédricos e tetraédricos; íons de metais de transição em sistemas biológicos5. Sistemas biológicos de transporte: Transporte de O2 em mamíferos, transferência de elétrons dependente de metaloproteínas;6. Processos catalíticos - ácido/base e oxido-redução em metaloproteínas: Proteínas contendo íon Zn2+, peroxidases, oxidases."

# 5) Row 19 (Metodo:) - replace B/C with evaluation method text
$ws.Range("B19").Value = "A avaliação será feita por meio de provas escritas."
$ws.Range("C19").Value = "A avaliação será feita por meio de provas escritas."

# 6) Row 20 (Criterio:) - replace B/C with NF formula text
$ws.Range("B20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."
$ws.Range("C20").Value = "A Nota final (NF) será calculada da seguinte maneira:NF = (P1 + 2*P2)/3Sendo que para P2 a matéria será cumulativa do semestre."

# 7) Row 21 (Norma de recuperacao:) - replace B/C with recovery text
$ws.Range("B21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"
$ws.Range("C21").Value = "A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR)/2"

# 8) Row 22 (Bibliografia:) - replace B/C with bibliography text
$ws.Range("B22").Value = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 20112. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"
$ws.Range("C22").Value = "1. Atkins e Jones, Princípios de Química, 5a edição, Bookman, 20112. Shiver e Atikins, Química Inorgânica, 4a edição, Bookman, 2008"

# 9) Split column A from the combined col(1:2) width group (cosmetic OOXML normalization)
$ws.Columns(1).ColumnWidth = 29.83

